# Commit: "@llmorreale will be lead author"
#
# The author list currently reads (in order):
#   Krystal Bagnaschi, Rachel Hoffman, Iris Kennedy, Erin MacMonigle,
#   Caroline Troy, Cameron Dow, Jen Jordan, Valentine Herrmann,
#   David Mitre, Jess Shue, Luca Morreale, William McShea, ...
#
# Luca Morreale needs to become the lead (first) author, i.e. his
# paragraph must move from just after "Jess Shue" to just before
# "Krystal Bagnaschi" (the current first author), leaving everything
# else untouched.

$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

# --- Step 1: insert the new lead-author paragraph right before the
#     current first author, "Krystal Bagnaschi". ---
$firstAuthor = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "Krystal Bagnaschi") {
        $firstAuthor = $p
        break
    }
}

if ($firstAuthor -eq $null) {
    throw "Could not find the 'Krystal Bagnaschi' paragraph."
}

# Insert text + paragraph mark immediately before the found paragraph;
# the new paragraph inherits the "BodyText" style of the paragraph it
# was inserted in front of.
$firstAuthor.Range.InsertBefore("Luca Morreale" + [char]13)

# --- Step 2: remove the old "Luca Morreale" paragraph, which sits
#     right after "Jess Shue". ---
$oldEntry = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "Luca Morreale") {
        $prev = $p.Previous()
        if ($prev -ne $null -and (Get-ParaText $prev) -eq "Jess Shue") {
            $oldEntry = $p
            break
        }
    }
}

if ($oldEntry -eq $null) {
    throw "Could not find the old 'Luca Morreale' paragraph after 'Jess Shue'."
}

$oldEntry.Range.Delete()
